# New attendance check-ins appended to the bottom of the log (rows 24-25).
# Column A (ID) is stored as text in the source data, so force a text
# number-format before assigning the numeric-looking id so Excel doesn't
# reinterpret it as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "1446896"
$ws.Cells.Item(24, 2).Value = "Unknown"
$ws.Cells.Item(24, 3).Value = "2025-01-10 00:49:41"

$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "1446896"
$ws.Cells.Item(25, 2).Value = "Unknown"
$ws.Cells.Item(25, 3).Value = "2025-01-10 00:52:00"
